# Change the fixed "date" footer placeholder text from 2019-04-26 to
# 01.05.2019 everywhere it is defined: the slide master and every slide
# layout (mirrors using the Header/Footer "Apply to All" flow in
# PowerPoint, which stores the cached date text on the master + layouts
# rather than on individual slides).

$p = $ppt.ActivePresentation
$newDate = "01.05.2019"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}
